$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert a brand-new data row at position 2 (pushes the existing
# rows 2-21 down to 3-22). Excel auto-copies the header row's formatting
# onto the inserted row, so strip that back off to match the plain data rows.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:H2").ClearFormats()

# Step 2: populate the newly inserted row 2 with its sensor reading.
$ws.Range("B2").Value = "walkingToRunning"
$ws.Range("C2").Value = -1.568910002708435
$ws.Range("D2").Value = -8.732925415039062
$ws.Range("E2").Value = 11.35853481292725
$ws.Range("F2").Value = 0.7064247653999316
$ws.Range("G2").Value = 0.3879705256011861
$ws.Range("H2").Value = -1.276251717872426

# Step 3: append 9 new trailing rows (rows 23-31) with fresh sensor readings,
# continuing on right after the shifted block, which now ends at row 22.
$ws.Range("B23").Value = "walkingToRunning"
$ws.Range("C23").Value = 18.49324798583984
$ws.Range("D23").Value = -80.07238006591797
$ws.Range("E23").Value = 68.10049438476562
$ws.Range("F23").Value = 3.552849229822288
$ws.Range("G23").Value = -5.358195619534956
$ws.Range("H23").Value = -0.1442216574237578

$ws.Range("B24").Value = "walkingToRunning"
$ws.Range("C24").Value = -9.736778259277344
$ws.Range("D24").Value = 7.034156322479248
$ws.Range("E24").Value = 9.938852310180664
$ws.Range("F24").Value = 4.79510967652812
$ws.Range("G24").Value = 16.44486069921317
$ws.Range("H24").Value = -4.148945381193625

$ws.Range("B25").Value = "walkingToRunning"
$ws.Range("C25").Value = 53.70849227905273
$ws.Range("D25").Value = -16.30745315551758
$ws.Range("E25").Value = 39.18264770507812
$ws.Range("F25").Value = -0.9187659382214614
$ws.Range("G25").Value = -0.2956210344576297
$ws.Range("H25").Value = 3.151463126168018

$ws.Range("B26").Value = "walkingToRunning"
$ws.Range("C26").Value = -80.71715545654297
$ws.Range("D26").Value = -15.45442008972168
$ws.Range("E26").Value = -28.87823677062988
$ws.Range("F26").Value = -5.69272972968636
$ws.Range("G26").Value = 3.907055351334684
$ws.Range("H26").Value = 4.745034019354134

$ws.Range("B27").Value = "walkingToRunning"
$ws.Range("C27").Value = 23.98992347717285
$ws.Range("D27").Value = 5.58967399597168
$ws.Range("E27").Value = -7.651249885559082
$ws.Range("F27").Value = -0.1373755151245123
$ws.Range("G27").Value = 11.3405332758948
$ws.Range("H27").Value = -3.450026544822681

$ws.Range("B28").Value = "walkingToRunning"
$ws.Range("C28").Value = 0.3348398208618164
$ws.Range("D28").Value = -21.52296257019043
$ws.Range("E28").Value = 16.69417762756348
$ws.Range("F28").Value = -1.20526529810765
$ws.Range("G28").Value = -4.696971940510186
$ws.Range("H28").Value = 1.208450563062909

$ws.Range("B29").Value = "walkingToRunning"
$ws.Range("C29").Value = 3.726076126098633
$ws.Range("D29").Value = 2.594820261001587
$ws.Range("E29").Value = 41.84358978271485
$ws.Range("F29").Value = 8.444541233142751
$ws.Range("G29").Value = -2.498347297840359
$ws.Range("H29").Value = -2.655560967885902

$ws.Range("B30").Value = "walkingToRunning"
$ws.Range("C30").Value = 8.552176475524902
$ws.Range("D30").Value = 7.943446636199951
$ws.Range("E30").Value = 19.50382232666016
$ws.Range("F30").Value = -0.723370986541445
$ws.Range("G30").Value = 0.8442579066087497
$ws.Range("H30").Value = -2.731481316125937

$ws.Range("B31").Value = "walkingToRunning"
$ws.Range("C31").Value = 28.44747161865234
$ws.Range("D31").Value = -58.03325653076172
$ws.Range("E31").Value = 40.07803344726562
$ws.Range("F31").Value = -3.363093618511545
$ws.Range("G31").Value = 1.143067340257809
$ws.Range("H31").Value = 1.45013582887988

# Step 4: the "timestamp" column A is simply 100 * (row-2) for every data row
# (row 2 -> 0, row 3 -> 100, ... row 31 -> 2900). Rewrite it explicitly for
# every data row so it is correct regardless of how the row-insert shifted
# the previously-existing cells.
For ($r = 2; $r -le 31; $r++) {
    $ws.Range("A" + $r).Value = ($r - 2) * 100
}